$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5431
$ws1.Range("F3").Value = 593
$ws1.Range("F4").Value = 11685
$ws1.Range("G4").Value = 62
$ws1.Range("F5").Value = 285
$ws1.Range("F6").Value = 594
$ws1.Range("F7").Value = 171
$ws1.Range("F8").Value = 280
$ws1.Range("F9").Value = 1031

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 24

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5431
$ws4.Range("F5").Value = 593
$ws4.Range("F6").Value = 24
$ws4.Range("F7").Value = 11685
$ws4.Range("G7").Value = 62
$ws4.Range("F8").Value = 285
$ws4.Range("F9").Value = 594
$ws4.Range("F10").Value = 171
$ws4.Range("F13").Value = 280
$ws4.Range("F14").Value = 1031
